$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("T13").Value = 15.579750353475299
$ws.Range("U13").Value = 17.021481171950001
$ws.Range("V13").Value = 17.021607898223198
$ws.Range("W13").Value = 0.0021543857891664201
$ws.Range("T14").Value = 15.6132855214661
$ws.Range("U14").Value = 17.900906098794401
$ws.Range("V14").Value = 17.8994459647305
$ws.Range("W14").Value = 0.0032601209394831098
$ws.Range("T15").Value = 15.682967789253899
$ws.Range("U15").Value = 18.894285164794901
$ws.Range("V15").Value = 18.8964740265816
$ws.Range("W15").Value = 0.0041889498838416899
$ws.Range("T16").Value = 15.3373289265569
$ws.Range("U16").Value = 1.87038660857343
$ws.Range("V16").Value = 1.8630447294889201
$ws.Range("W16").Value = 0.000021448621537529901
$ws.Range("T17").Value = 15.3912312075096
$ws.Range("U17").Value = 3.39060418704676
$ws.Range("V17").Value = 3.3904720877499801
$ws.Range("W17").Value = 0.00032781251260287001
$ws.Range("T18").Value = 15.4012446651919
$ws.Range("U18").Value = 3.6888090138851899
$ws.Range("V18").Value = 3.6890717142934499
$ws.Range("W18").Value = 0.00097116606816001299
$ws.Range("T19").Value = 15.7368213871005
$ws.Range("U19").Value = 1.17716779061782
$ws.Range("V19").Value = 1.17716779061782
$ws.Range("W19").Value = 0.00000031342752169547099
$ws.Range("T20").Value = 15.334430922831899
$ws.Range("U20").Value = 0.59541791989840198
$ws.Range("V20").Value = 0.59541407266858204
$ws.Range("W20").Value = 0.000061327053976945896
$ws.Range("T21").Value = 15.3748806484939
$ws.Range("U21").Value = 0.86675747048060103
$ws.Range("V21").Value = 0.86691668342397898
$ws.Range("W21").Value = 0.00046926626837150799
$ws.Range("T22").Value = 14.516634211657401
$ws.Range("U22").Value = 13.0758831169444
$ws.Range("V22").Value = 13.076827957316601
$ws.Range("W22").Value = 0.00161864081803668
$ws.Range("T23").Value = 14.539535515907801
$ws.Range("U23").Value = 15.245760178260699
$ws.Range("V23").Value = 15.227502768254899
$ws.Range("W23").Value = 0.0025803986264241502
$ws.Range("T24").Value = 14.550231060266899
$ws.Range("U24").Value = 15.5013665940139
$ws.Range("V24").Value = 15.468692017594201
$ws.Range("W24").Value = 0.0034803908842232201
$ws.Range("T25").Value = 13.786223762814
$ws.Range("U25").Value = 0.66063657160254396
$ws.Range("V25").Value = 0.68884895646424604
$ws.Range("W25").Value = 0.000016852065228841999
$ws.Range("T26").Value = 14.0070021225304
$ws.Range("U26").Value = 1.06678369869075
$ws.Range("V26").Value = 1.04325017432579
$ws.Range("W26").Value = 0.000172901989445172
$ws.Range("T27").Value = 14.1488570309546
$ws.Range("U27").Value = 2.1309650170560501
$ws.Range("V27").Value = 2.1176026617924699
$ws.Range("W27").Value = 0.00068231060381754999
$ws.Range("T28").Value = 14.4019086649791
$ws.Range("U28").Value = 0.91874436674812898
$ws.Range("V28").Value = 1.00411458314247
$ws.Range("W28").Value = 0.000000313427521695926
$ws.Range("T29").Value = 13.9437087813871
$ws.Range("U29").Value = 0.115517870908374
$ws.Range("V29").Value = 0.115517870908374
$ws.Range("W29").Value = 0.000054178941282424999
$ws.Range("T30").Value = 13.6603035280811
$ws.Range("U30").Value = 0.024444814902255401
$ws.Range("V30").Value = 0.0303957250212098
$ws.Range("W30").Value = 0.00042590482568528301
$ws.Range("T31").Value = 13.4154865178308
$ws.Range("U31").Value = 8.8838905219535906
$ws.Range("V31").Value = 8.7815108992329005
$ws.Range("W31").Value = 0.0014996423274823401
$ws.Range("T32").Value = 13.488212088292601
$ws.Range("U32").Value = 11.6011871992241
$ws.Range("V32").Value = 11.464277009181099
$ws.Range("W32").Value = 0.0023055782250070902
$ws.Range("T33").Value = 13.4493877525701
$ws.Range("U33").Value = 11.254330159209401
$ws.Range("V33").Value = 11.1644269742297
$ws.Range("W33").Value = 0.0031948337357229901
$ws.Range("T34").Value = 12.7074251381341
$ws.Range("U34").Value = 0.223412643173623
$ws.Range("V34").Value = 0.26917716361289501
$ws.Range("W34").Value = 0.000016850155423055
$ws.Range("T35").Value = 12.655349256227399
$ws.Range("U35").Value = 0.181339652993429
$ws.Range("V35").Value = 0.18362453841707499
$ws.Range("W35").Value = 0.00016539494629814699
$ws.Range("T36").Value = 12.9790954608777
$ws.Range("U36").Value = 0.43923130490677698
$ws.Range("V36").Value = 0.46084613765128701
$ws.Range("W36").Value = 0.000619578282181649
$ws.Range("T37").Value = 13.257740588890201
$ws.Range("U37").Value = 0.46546668697346799
$ws.Range("V37").Value = 0.46664741860531
$ws.Range("W37").Value = 0.00000031364705384216898
$ws.Range("T38").Value = 12.454552110318801
$ws.Range("U38").Value = 0.012684090681526801
$ws.Range("V38").Value = 0.0125866161755198
$ws.Range("W38").Value = 0.000054463020771851297
$ws.Range("T39").Value = 12.4923832433483
$ws.Range("U39").Value = 0.013474105329702999
$ws.Range("V39").Value = 0.012552358726014101
$ws.Range("W39").Value = 0.00042591225612303302
$ws.Range("T40").Value = 11.910637107833301
$ws.Range("U40").Value = 1.29224480095224
$ws.Range("V40").Value = 1.2324058648476
$ws.Range("W40").Value = 0.0014996667658895501
$ws.Range("T41").Value = 11.957051337013301
$ws.Range("U41").Value = 1.57724608827539
$ws.Range("V41").Value = 1.6122134527930001
$ws.Range("W41").Value = 0.0023000794293604201
$ws.Range("T42").Value = 11.919627832740201
$ws.Range("U42").Value = 1.5161627818275401
$ws.Range("V42").Value = 1.5089374337874899
$ws.Range("W42").Value = 0.0031837894900912799
$ws.Range("T43").Value = 11.763716808226199
$ws.Range("U43").Value = 0.087305542854317103
$ws.Range("V43").Value = 0.091078086438149602
$ws.Range("W43").Value = 0.000016849893719218501
$ws.Range("T44").Value = 11.749854916832
$ws.Range("U44").Value = 0.082241369587612204
$ws.Range("V44").Value = 0.079585171189683196
$ws.Range("W44").Value = 0.000164344451552019
$ws.Range("T45").Value = 11.782094577052201
$ws.Range("U45").Value = 0.093722710975636905
$ws.Range("V45").Value = 0.094752640655732398
$ws.Range("W45").Value = 0.00061958095290298297
$ws.Range("T46").Value = 12.977945701256299
$ws.Range("U46").Value = 0.99823784259652104
$ws.Range("V46").Value = 0.99823784259652104
$ws.Range("W46").Value = 0.00000031342752169547099
$ws.Range("T47").Value = 12.4393573218436
$ws.Range("U47").Value = 0.0322329954563546
$ws.Range("V47").Value = 0.033050250417763299
$ws.Range("W47").Value = 0.000054195746531310103
$ws.Range("T48").Value = 11.77435193354
$ws.Range("U48").Value = 0.0138253393984578
$ws.Range("V48").Value = 0.0126810514797586
$ws.Range("W48").Value = 0.00042590077871337199

$ws.Activate()
$ws.Range("A2").Select()
